$d = $word.ActiveDocument

$section = $d.Sections.First
$footer = $section.Footers.Item(1)  # wdHeaderFooterPrimary = 1

$footer.PageNumbers.Add(2, $false)  # wdAlignPageNumberRight = 2, firstPage = $false
